$d = $word.ActiveDocument

$pairs = @(
    @("53÷2=26, 1", "28÷4=7, 0"),
    @("96÷7=13, 5", "84÷8=10, 4"),
    @("22÷3=7, 1", "43÷9=4, 7"),
    @("32÷8=4, 0", "92÷3=30, 2"),
    @("77÷5=15, 2", "84÷5=16, 4"),
    @("41÷8=5, 1", "87÷8=10, 7"),
    @("49÷7=7, 0", "28÷6=4, 4"),
    @("37÷5=7, 2", "75÷3=25, 0"),
    @("41÷6=6, 5", "89÷3=29, 2"),
    @("71÷8=8, 7", "54÷9=6, 0"),
    @("56÷8=7, 0", "60÷2=30, 0"),
    @("95÷5=19, 0", "99÷8=12, 3"),
    @("47÷5=9, 2", "49÷2=24, 1"),
    @("22÷2=11, 0", "35÷4=8, 3"),
    @("11÷4=2, 3", "61÷9=6, 7"),
    @("23÷3=7, 2", "37÷2=18, 1"),
    @("43÷6=7, 1", "17÷2=8, 1"),
    @("89÷8=11, 1", "63÷4=15, 3"),
    @("41÷7=5, 6", "30÷9=3, 3"),
    @("26÷8=3, 2", "81÷5=16, 1"),
    @("85÷5=17, 0", "57÷6=9, 3"),
    @("21÷5=4, 1", "92÷3=30, 2"),
    @("32÷3=10, 2", "75÷8=9, 3"),
    @("39÷3=13, 0", "92÷3=30, 2"),
    @("16÷4=4, 0", "73÷8=9, 1")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
